$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 10000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -10350
$ws.Range("H64").Value = 3401.5
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3401.5
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3401.5
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3897.5
$ws.Range("H67").Value = 3401.5
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3401.5
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3401.5
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -5117.5
$ws.Range("H86").Value = 7850
$ws.Range("I86").Value = 9000
$ws.Range("J86").Value = 7466.6665
$ws.Range("K86").Value = 9000
$ws.Range("L86").Value = 7466.6665
$ws.Range("M86").Value = -7877
$ws.Range("N86").Value = -9712.666499999999
$ws.Range("H89").Value = 7850
$ws.Range("I89").Value = 9000
$ws.Range("J89").Value = 7466.6665
$ws.Range("K89").Value = 45000
$ws.Range("L89").Value = 37333.3325
$ws.Range("M89").Value = -39384
$ws.Range("N89").Value = -48565.3325

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 799.5
$ws.Range("I2").Value = 799.5
$ws.Range("K2").Value = 799.5
$ws.Range("M2").Value = -686.5
$ws.Range("H32").Value = 4813.6294
$ws.Range("I32").Value = 3413.524
$ws.Range("J32").Value = 9714
$ws.Range("K32").Value = 3413.524
$ws.Range("L32").Value = 9714
$ws.Range("M32").Value = -3126.524
$ws.Range("N32").Value = -10288
$ws.Range("H116").Value = 799.5
$ws.Range("I116").Value = 799.5
$ws.Range("K116").Value = 799.5
$ws.Range("M116").Value = 1494.5
$ws.Range("H122").Value = 1764.52
$ws.Range("I122").Value = 1532.5264
$ws.Range("K122").Value = 4597.5792
$ws.Range("M122").Value = -2147.5792
$ws.Range("H132").Value = 1750
$ws.Range("I132").Value = 1750
$ws.Range("K132").Value = 5250
$ws.Range("M132").Value = -2720

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 799.5
$ws.Range("I3").Value = 799.5
$ws.Range("K3").Value = 799.5
$ws.Range("M3").Value = -685.5
$ws.Range("H94").Value = 2667.375
$ws.Range("I94").Value = 2334.1428
$ws.Range("K94").Value = 2334.1428
$ws.Range("M94").Value = -1883.1428
$ws.Range("H134").Value = 712
$ws.Range("I134").Value = 712
$ws.Range("K134").Value = 2136
$ws.Range("M134").Value = 399

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3707.0625
$ws.Range("I132").Value = 3371.077
$ws.Range("J132").Value = 5163
$ws.Range("K132").Value = 10113.231
$ws.Range("L132").Value = 15489
$ws.Range("M132").Value = -7583.231
$ws.Range("N132").Value = -20549

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 250581
$ws.Range("I4").Value = 775
$ws.Range("J4").Value = 999999
$ws.Range("K4").Value = 2325
$ws.Range("L4").Value = 2999997
$ws.Range("M4").Value = -2213
$ws.Range("N4").Value = -3000221
$ws.Range("H132").Value = 2500.5715
$ws.Range("I132").Value = 2500.8
$ws.Range("K132").Value = 22507.2
$ws.Range("M132").Value = -19977.2
$ws.Range("H133").Value = 5000
$ws.Range("J133").Value = 5000
$ws.Range("L133").Value = 15000
$ws.Range("N133").Value = -25120
$ws.Range("H134").Value = 2000
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -930
$ws.Range("H138").Value = 3000
$ws.Range("I138").Value = 3000
$ws.Range("K138").Value = 9000
$ws.Range("M138").Value = -3860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2634.8333
$ws.Range("I80").Value = 2003
$ws.Range("J80").Value = 3266.6667
$ws.Range("K80").Value = 2003
$ws.Range("L80").Value = 3266.6667
$ws.Range("M80").Value = -1005
$ws.Range("N80").Value = -5262.6667
$ws.Range("H83").Value = 2634.8333
$ws.Range("I83").Value = 2003
$ws.Range("J83").Value = 3266.6667
$ws.Range("K83").Value = 10015
$ws.Range("L83").Value = 16333.3335
$ws.Range("M83").Value = -5023
$ws.Range("N83").Value = -26317.3335
$ws.Range("H102").Value = 3480.1333
$ws.Range("I102").Value = 3609.4546
$ws.Range("J102").Value = 3124.5
$ws.Range("K102").Value = 3609.4546
$ws.Range("L102").Value = 3124.5
$ws.Range("M102").Value = -1987.4546
$ws.Range("N102").Value = -6368.5
$ws.Range("H132").Value = 3092.375
$ws.Range("I132").Value = 3092.375
$ws.Range("K132").Value = 9277.125
$ws.Range("M132").Value = -6747.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3611.3333
$ws.Range("I61").Value = 3757.9092
$ws.Range("J61").Value = 1999
$ws.Range("K61").Value = 3757.9092
$ws.Range("L61").Value = 1999
$ws.Range("M61").Value = -3555.9092
$ws.Range("N61").Value = -2403
$ws.Range("H82").Value = 2198.077
$ws.Range("I82").Value = 1999.5
$ws.Range("J82").Value = 2515.8
$ws.Range("K82").Value = 1999.5
$ws.Range("L82").Value = 2515.8
$ws.Range("M82").Value = -1638.5
$ws.Range("N82").Value = -3237.8
$ws.Range("H85").Value = 2198.077
$ws.Range("I85").Value = 1999.5
$ws.Range("J85").Value = 2515.8
$ws.Range("K85").Value = 1999.5
$ws.Range("L85").Value = 2515.8
$ws.Range("M85").Value = -751.5
$ws.Range("N85").Value = -5011.8
$ws.Range("H93").Value = 5000
$ws.Range("I93").Value = 5000
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 5000
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -3752
$ws.Range("N93").ClearContents()
$ws.Range("H113").Value = 3611.3333
$ws.Range("I113").Value = 3757.9092
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 3757.9092
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = -1587.9092
$ws.Range("N113").Value = -6339
$ws.Range("H122").Value = 3482.7144
$ws.Range("I122").Value = 3721.7693
$ws.Range("J122").Value = 375
$ws.Range("K122").Value = 11165.3079
$ws.Range("L122").Value = 1125
$ws.Range("M122").Value = -8715.3079
$ws.Range("N122").Value = -6025

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 49943.5
$ws.Range("J74").Value = 49943.5
$ws.Range("L74").Value = 49943.5
$ws.Range("N74").Value = -51815.5
$ws.Range("H77").Value = 49943.5
$ws.Range("J77").Value = 49943.5
$ws.Range("L77").Value = 149830.5
$ws.Range("N77").Value = -159190.5
$ws.Range("H107").Value = 543.5454999999999
$ws.Range("I107").Value = 453.22223
$ws.Range("K107").Value = 1359.66669
$ws.Range("M107").Value = 560.33331
$ws.Range("H113").Value = 824.25
$ws.Range("I113").Value = 765.6667
$ws.Range("K113").Value = 2297.0001
$ws.Range("M113").Value = -127.0001000000002
$ws.Range("H122").Value = 167707.58
$ws.Range("I122").Value = 167707.58
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 503122.74
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -500672.74
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3237.182
$ws.Range("I126").Value = 2656.6667
$ws.Range("J126").Value = 5849.5
$ws.Range("K126").Value = 7970.000100000001
$ws.Range("L126").Value = 17548.5
$ws.Range("M126").Value = -5500.000100000001
$ws.Range("N126").Value = -22488.5
